$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.619.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.563.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'210.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "'24.82"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Value = "'0.0587"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Value = "'1.787.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'1.560.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'28.656.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Value = "'61.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'227.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").Value = "'3.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'9.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'151.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'14.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("D33").Value = "'1.401.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("D39").Value = "'0.0163"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.517"
$ws.Range("D40").Style = "Normal"
$ws.Range("D45").Value = "'63.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'1.699.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.843"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'84.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'42.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.0512"
$ws.Range("D51").Style = "Normal"

$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("E6").Value = "  +6.11%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +3.58%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("E35").Value = "  -3.54%  "
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("E48").Value = "  -8.45%  "
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("E50").Value = "  +4.44%  "
$ws.Range("E51").Value = "  -0.49%  "
